# Update cryptos list (Price / Volume(1h) columns, and the two swapped
# Coin/Link rows 32-33) to match the latest scraped values.
#
# NumberFormat "@" (Text) is applied before setting each value so that
# Excel keeps numeric-looking strings (e.g. "1.00", "364.72") as literal
# text instead of silently coercing them to numbers; Style is reset back
# to "Normal" right after so the cell's on-disk style stays the same
# (unstyled) as before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.068.85'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.44%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.910.39'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '364.72'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.25'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.41%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.540'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.91%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.589'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -7.22%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.95'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.57%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.41'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.17%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.374.74'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.58%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.35'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -6.11%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.915.59'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.950'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.78%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.081.32'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.35%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.31'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -7.28%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.25'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.35%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.02'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -6.69%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0946'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.93%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.09'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.63%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.55'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.40%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.68'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.91%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.77%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.96%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.90'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.34%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.29'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.21%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.104'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("B32").Style = "Normal"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C32").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.94'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.59%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("B33").Style = "Normal"

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C33").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.10'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '35.11'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.26%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '50.66'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.22%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0422'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.96%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.81'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.13'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.58%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.93'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -7.61%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.86'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.14%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.50%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.56'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '117.70'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.03%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.10'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.066.70'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.11%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.20'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -7.93%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -8.45%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.205.56'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.87%  '
$ws.Range("E51").Style = "Normal"
